# Apply the "Softexpert" update:
#  - Organograma: insert 2 new rows after row 18 (duplicating the existing
#    row 18 "Key User Documentos / Ana Martins / Gestora de Qualidade / ..."
#    entry) and change column E for rows 18-20 from "Qualidade & Compliance"
#    to "TI". Update autofilter / dimension / view selections accordingly.
#  - Projetos: move the frozen-pane anchor to D13.
#  - Riscos: set the active selection to C17.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Organograma sheet
# ---------------------------------------------------------------------
$org = $wb.Worksheets.Item("Organograma")
$org.Activate()

# Duplicate row 18 twice by copying it and inserting the copied cells as new
# rows immediately below - this carries the source row's formatting
# (borders, fonts, row height) along with it.
$org.Rows.Item(18).Copy()
$org.Rows.Item(19).Insert()
$org.Rows.Item(18).Copy()
$org.Rows.Item(19).Insert()
$excel.CutCopyMode = $false

# The row-insert above can lose the per-cell border formatting of the
# source row on the newly-created rows; re-apply row 18's formatting
# (borders/fonts/fill/number-format) explicitly onto rows 19:20.
$org.Range("A18:E18").Copy()
$org.Range("A19:E20").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Update column E for the three duplicated rows to "TI".
$org.Range("E18:E20").Value = "TI"

# These three rows render at 15pt (matching the sheet's other single-line
# data rows) once they pick up the row-18 formatting.
$org.Range("A18:E20").RowHeight = 15

# Refresh the autofilter range to match the new data extent. The sheet
# already has an active AutoFilter, so toggle it off then back on against
# the new A1:E33 extent (a single call would just switch the existing one
# off).
$org.Range("A1:E33").AutoFilter()
$org.Range("A1:E33").AutoFilter()

# Keep the workbook-level hidden _FilterDatabase name in sync with the new
# autofilter extent.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
  $n = $wb.Names.Item($i)
  if ($n.Name -eq "Organograma!_FilterDatabase") {
    $n.RefersTo = "=Organograma!`$A`$1:`$E`$33"
  }
}

# Restore the selection: A20:E20 (the header rows stay frozen at ySplit=2).
$org.Range("A20:E20").Select()

# ---------------------------------------------------------------------
# Projetos sheet
# ---------------------------------------------------------------------
$proj = $wb.Worksheets.Item("Projetos")
$proj.Activate()
$proj.Range("D13").Select()

# ---------------------------------------------------------------------
# Riscos sheet
# ---------------------------------------------------------------------
$risc = $wb.Worksheets.Item("Riscos")
$risc.Activate()
$risc.Range("C17").Select()

# Leave Organograma as the active sheet (matches tabSelected in the sheet XML).
$org.Activate()
